# Rename the document title:
#  - remove the leading empty "Corpodetexto" paragraph
#  - remove the "DECLARAÇÃO DE VISITA" Título paragraph
#  - turn the "(Arts 722 e 729 ...)" paragraph into a bold
#    "AVALIAÇÃO DE EXPERIÊNCIA DO CLIENTE" paragraph
$d = $word.ActiveDocument

# Paragraph 2 is the "Ttulo" paragraph containing "DECLARAÇÃO DE VISITA".
# Deleting its own Range (which includes its own paragraph mark) removes
# that paragraph node outright, without disturbing its neighbours.
$titlePara = $d.Paragraphs.Item(2)
$titlePara.Range.Delete()

# Paragraph 1 is now the empty "Corpodetexto" spacer paragraph above the
# title. Deleting its Range removes that paragraph node too, leaving the
# former "(Arts ...)" paragraph as the new first paragraph, keeping that
# paragraph's own identity/formatting (pPr).
$emptyPara = $d.Paragraphs.Item(1)
$emptyPara.Range.Delete()

# This is now the "(Arts 722 e 729 do Código Civil c/c Art. 20 da Lei
# n°6.530/78)" paragraph - replace its text with the new title and make
# it bold.
$p = $d.Paragraphs.Item(1)

$textRange = $p.Range
$textRange.End = $textRange.End - 1
$textRange.Text = "AVALIAÇÃO DE EXPERIÊNCIA DO CLIENTE"

# Apply bold (regular + complex-script) to the whole paragraph, including
# its paragraph mark, so both the run and the paragraph formatting pick
# up the bold attribute.
$fullRange = $p.Range
$fullRange.Bold = 1
$fullRange.Font.BoldBi = 1
